{"js": "// Replace each arithmetic equation's answer text in the worksheet table\n// with its new value. Each entry is a unique, literal old->new text pair\n// derived from the target diff; only the <w:t> run text changes, so we\n// search for the exact old text and do an in-place text replace (this\n// preserves the existing run/paragraph formatting: fonts, size, alignment).\nconst replacements = [\n  [\"1+36=37\", \"44+24=68\"],\n  [\"79-0=79\", \"54-16=38\"],\n  [\"29-21=8\", \"13-7=6\"],\n  [\"81-48=33\", \"52+35=87\"],\n  [\"34+26=60\", \"99-88=11\"],\n  [\"53-29=24\", \"53-46=7\"],\n  [\"37+61=98\", \"81-2=79\"],\n  [\"96-54=42\", \"74+10=84\"],\n  [\"7+46=53\", \"53+14=67\"],\n  [\"51-26=25\", \"30+4=34\"],\n  [\"85-57=28\", \"30+14=44\"],\n  [\"5+35=40\", \"9+61=70\"],\n  [\"24+29=53\", \"60-35=25\"],\n  [\"42-29=13\", \"54+11=65\"],\n  [\"23+40=63\", \"69-49=20\"],\n  [\"87-36=51\", \"77+13=90\"],\n  [\"23+22=45\", \"18-0=18\"],\n  [\"29+41=70\", \"82-58=24\"],\n  [\"82-64=18\", \"35-22=13\"],\n  [\"32-5=27\", \"83+11=94\"],\n  [\"69+28=97\", \"39+9=48\"],\n  [\"13-3=10\", \"12+69=81\"],\n  [\"59+39=98\", \"56-27=29\"],\n  [\"56+31=87\", \"52+10=62\"],\n  [\"76-33=43\", \"0+32=32\"],\n  [\"95-69=26\", \"57-37=20\"],\n  [\"6+66=72\", \"86-10=76\"],\n  [\"36-8=28\", \"78-3=75\"],\n  [\"59-18=41\", \"13+4=17\"],\n  [\"61-44=17\", \"70-65=5\"],\n  [\"81-23=58\", \"97-93=4\"],\n  [\"31+4=35\", \"63-34=29\"],\n  [\"56-16=40\", \"7+73=80\"],\n  [\"64+3=67\", \"61-32=29\"],\n  [\"83-21=62\", \"62-48=14\"],\n  [\"2+22=24\", \"75-48=27\"],\n  [\"35+22=57\", \"70+3=73\"],\n  [\"31-20=11\", \"99-33=66\"],\n  [\"78-76=2\", \"44+22=66\"],\n  [\"86+10=96\", \"68-19=49\"],\n  [\"96-91=5\", \"51-38=13\"],\n  [\"67-61=6\", \"32-32=0\"],\n  [\"43+22=65\", \"21-16=5\"],\n  [\"40-23=17\", \"96-51=45\"],\n  [\"21+45=66\", \"43-23=20\"],\n  [\"18+71=89\", \"78-35=43\"],\n  [\"58-0=58\", \"75-61=14\"],\n  [\"66+16=82\", \"25+37=62\"],\n  [\"10+3=13\", \"79-63=16\"],\n  [\"44+28=72\", \"25+72=97\"],\n  [\"42-35=7\", \"4+74=78\"],\n  [\"57-27=30\", \"72+24=96\"],\n  [\"68-0=68\", \"0+20=20\"],\n  [\"76-25=51\", \"2+44=46\"],\n  [\"68-66=2\", \"98-46=52\"],\n  [\"89-20=69\", \"67-8=59\"],\n  [\"29+6=35\", \"26+31=57\"],\n  [\"53-49=4\", \"40+52=92\"],\n  [\"17+76=93\", \"84-0=84\"],\n  [\"59+1=60\", \"41-27=14\"],\n  [\"28+30=58\", \"32+66=98\"],\n  [\"18-2=16\", \"66+8=74\"],\n  [\"0+3=3\", \"62+19=81\"],\n  [\"16+17=33\", \"9-7=2\"],\n  [\"28-12=16\", \"99-67=32\"],\n  [\"42+40=82\", \"26+20=46\"],\n  [\"5+14=19\", \"61+21=82\"],\n  [\"55-45=10\", \"73+14=87\"],\n  [\"11+40=51\", \"64-41=23\"],\n  [\"44+48=92\", \"53-46=7\"],\n  [\"69-29=40\", \"7+67=74\"],\n  [\"59-47=12\", \"99-5=94\"],\n  [\"71+7=78\", \"91-21=70\"],\n  [\"98-76=22\", \"64-23=41\"],\n  [\"78-47=31\", \"77-62=15\"],\n  [\"38+19=57\", \"28+18=46\"],\n  [\"17+63=80\", \"3+6=9\"],\n  [\"75+23=98\", \"64-4=60\"],\n  [\"43-9=34\", \"88-16=72\"],\n  [\"53-7=46\", \"85-46=39\"],\n  [\"50+28=78\", \"17+13=30\"],\n  [\"67-62=5\", \"26+47=73\"],\n  [\"82-76=6\", \"6+47=53\"],\n  [\"6+69=75\", \"39-10=29\"],\n  [\"45-17=28\", \"98-30=68\"],\n  [\"85-52=33\", \"53+24=77\"],\n  [\"74+4=78\", \"91-54=37\"],\n  [\"36-27=9\", \"80-29=51\"],\n  [\"35-20=15\", \"25-13=12\"],\n  [\"71+4=75\", \"46+46=92\"],\n  [\"7+75=82\", \"46+44=90\"],\n  [\"38+49=87\", \"92-91=1\"],\n  [\"89-63=26\", \"80-59=21\"],\n  [\"79-51=28\", \"63-12=51\"],\n  [\"14+10=24\", \"89-87=2\"],\n  [\"80-73=7\", \"24+24=48\"],\n  [\"83-49=34\", \"98-71=27\"],\n  [\"54+26=80\", \"67+16=83\"],\n  [\"46-11=35\", \"94-71=23\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n\n  // Replace every match in place (old values are unique in this document,\n  // so this is normally exactly one hit).\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the arithmetic-answer worksheet table: replace each equation's\n# answer text with its new value. Each pair is a unique, literal old->new\n# text value derived from the target diff; only the run text changes (the\n# existing run/paragraph formatting -- fonts, size, alignment -- is left\n# untouched by Find/Replace).\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"1+36=37\", \"44+24=68\"),\n    @(\"79-0=79\", \"54-16=38\"),\n    @(\"29-21=8\", \"13-7=6\"),\n    @(\"81-48=33\", \"52+35=87\"),\n    @(\"34+26=60\", \"99-88=11\"),\n    @(\"53-29=24\", \"53-46=7\"),\n    @(\"37+61=98\", \"81-2=79\"),\n    @(\"96-54=42\", \"74+10=84\"),\n    @(\"7+46=53\", \"53+14=67\"),\n    @(\"51-26=25\", \"30+4=34\"),\n    @(\"85-57=28\", \"30+14=44\"),\n    @(\"5+35=40\", \"9+61=70\"),\n    @(\"24+29=53\", \"60-35=25\"),\n    @(\"42-29=13\", \"54+11=65\"),\n    @(\"23+40=63\", \"69-49=20\"),\n    @(\"87-36=51\", \"77+13=90\"),\n    @(\"23+22=45\", \"18-0=18\"),\n    @(\"29+41=70\", \"82-58=24\"),\n    @(\"82-64=18\", \"35-22=13\"),\n    @(\"32-5=27\", \"83+11=94\"),\n    @(\"69+28=97\", \"39+9=48\"),\n    @(\"13-3=10\", \"12+69=81\"),\n    @(\"59+39=98\", \"56-27=29\"),\n    @(\"56+31=87\", \"52+10=62\"),\n    @(\"76-33=43\", \"0+32=32\"),\n    @(\"95-69=26\", \"57-37=20\"),\n    @(\"6+66=72\", \"86-10=76\"),\n    @(\"36-8=28\", \"78-3=75\"),\n    @(\"59-18=41\", \"13+4=17\"),\n    @(\"61-44=17\", \"70-65=5\"),\n    @(\"81-23=58\", \"97-93=4\"),\n    @(\"31+4=35\", \"63-34=29\"),\n    @(\"56-16=40\", \"7+73=80\"),\n    @(\"64+3=67\", \"61-32=29\"),\n    @(\"83-21=62\", \"62-48=14\"),\n    @(\"2+22=24\", \"75-48=27\"),\n    @(\"35+22=57\", \"70+3=73\"),\n    @(\"31-20=11\", \"99-33=66\"),\n    @(\"78-76=2\", \"44+22=66\"),\n    @(\"86+10=96\", \"68-19=49\"),\n    @(\"96-91=5\", \"51-38=13\"),\n    @(\"67-61=6\", \"32-32=0\"),\n    @(\"43+22=65\", \"21-16=5\"),\n    @(\"40-23=17\", \"96-51=45\"),\n    @(\"21+45=66\", \"43-23=20\"),\n    @(\"18+71=89\", \"78-35=43\"),\n    @(\"58-0=58\", \"75-61=14\"),\n    @(\"66+16=82\", \"25+37=62\"),\n    @(\"10+3=13\", \"79-63=16\"),\n    @(\"44+28=72\", \"25+72=97\"),\n    @(\"42-35=7\", \"4+74=78\"),\n    @(\"57-27=30\", \"72+24=96\"),\n    @(\"68-0=68\", \"0+20=20\"),\n    @(\"76-25=51\", \"2+44=46\"),\n    @(\"68-66=2\", \"98-46=52\"),\n    @(\"89-20=69\", \"67-8=59\"),\n    @(\"29+6=35\", \"26+31=57\"),\n    @(\"53-49=4\", \"40+52=92\"),\n    @(\"17+76=93\", \"84-0=84\"),\n    @(\"59+1=60\", \"41-27=14\"),\n    @(\"28+30=58\", \"32+66=98\"),\n    @(\"18-2=16\", \"66+8=74\"),\n    @(\"0+3=3\", \"62+19=81\"),\n    @(\"16+17=33\", \"9-7=2\"),\n    @(\"28-12=16\", \"99-67=32\"),\n    @(\"42+40=82\", \"26+20=46\"),\n    @(\"5+14=19\", \"61+21=82\"),\n    @(\"55-45=10\", \"73+14=87\"),\n    @(\"11+40=51\", \"64-41=23\"),\n    @(\"44+48=92\", \"53-46=7\"),\n    @(\"69-29=40\", \"7+67=74\"),\n    @(\"59-47=12\", \"99-5=94\"),\n    @(\"71+7=78\", \"91-21=70\"),\n    @(\"98-76=22\", \"64-23=41\"),\n    @(\"78-47=31\", \"77-62=15\"),\n    @(\"38+19=57\", \"28+18=46\"),\n    @(\"17+63=80\", \"3+6=9\"),\n    @(\"75+23=98\", \"64-4=60\"),\n    @(\"43-9=34\", \"88-16=72\"),\n    @(\"53-7=46\", \"85-46=39\"),\n    @(\"50+28=78\", \"17+13=30\"),\n    @(\"67-62=5\", \"26+47=73\"),\n    @(\"82-76=6\", \"6+47=53\"),\n    @(\"6+69=75\", \"39-10=29\"),\n    @(\"45-17=28\", \"98-30=68\"),\n    @(\"85-52=33\", \"53+24=77\"),\n    @(\"74+4=78\", \"91-54=37\"),\n    @(\"36-27=9\", \"80-29=51\"),\n    @(\"35-20=15\", \"25-13=12\"),\n    @(\"71+4=75\", \"46+46=92\"),\n    @(\"7+75=82\", \"46+44=90\"),\n    @(\"38+49=87\", \"92-91=1\"),\n    @(\"89-63=26\", \"80-59=21\"),\n    @(\"79-51=28\", \"63-12=51\"),\n    @(\"14+10=24\", \"89-87=2\"),\n    @(\"80-73=7\", \"24+24=48\"),\n    @(\"83-49=34\", \"98-71=27\"),\n    @(\"54+26=80\", \"67+16=83\"),\n    @(\"46-11=35\", \"94-71=23\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n\n    # Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards,\n    #         MatchSoundsLike, MatchAllWordForms, Forward, Wrap, Format,\n    #         ReplaceWith, Replace)\n    $found = $find.Execute($oldText, $true, $true, $false, $false, $false, $true, 1, $false, $newText, 2)\n\n    if (-not $found) {\n        throw \"No match found for: $oldText\"\n    }\n}\n"}
